$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 506; existing rows 506-542 shift down to 507-543.
$ws.Rows.Item(506).Insert()

# Populate the newly inserted row 506 with the new weekly record.
$ws.Cells.Item(506, 1).Value  = 4
$ws.Cells.Item(506, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(506, 3).Value  = "Los Lagos"
$ws.Cells.Item(506, 4).Value  = 45265
$ws.Cells.Item(506, 5).Value  = 10
$ws.Cells.Item(506, 6).Value  = 100112043
$ws.Cells.Item(506, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(506, 8).Value  = "Sin especificar"
$ws.Cells.Item(506, 9).Value  = "Primera"
$ws.Cells.Item(506, 10).Value = 400
$ws.Cells.Item(506, 11).Value = 21000
$ws.Cells.Item(506, 12).Value = 21000
$ws.Cells.Item(506, 13).Value = 21000
$ws.Cells.Item(506, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(506, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(506, 16).Value = 350
$ws.Cells.Item(506, 17).Value = 60
$ws.Cells.Item(506, 18).Value = "Hortaliza"
